$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the leading "NN-0000  " occupation-code prefix from column B values,
# leaving just the occupation title text (column A / major names are unchanged).
$ws.Cells.Item(2, 2).Value  = "Farming, Fishing, and Forestry Occupations"
$ws.Cells.Item(3, 2).Value  = "Arts, Design, Entertainment, Sports, and Media Occupations"
$ws.Cells.Item(4, 2).Value  = "Arts, Design, Entertainment, Sports, and Media Occupations"
$ws.Cells.Item(5, 2).Value  = "Arts, Design, Entertainment, Sports, and Media Occupations"
$ws.Cells.Item(6, 2).Value  = "Arts, Design, Entertainment, Sports, and Media Occupations"
$ws.Cells.Item(7, 2).Value  = "Business and Financial Operations Occupations"
$ws.Cells.Item(8, 2).Value  = "Management Occupations"
$ws.Cells.Item(9, 2).Value  = "Computer and Mathematical Occupations"
$ws.Cells.Item(10, 2).Value = "Educational Instruction and Library Occupations"
$ws.Cells.Item(11, 2).Value = "Architecture and Engineering Occupations"
$ws.Cells.Item(12, 2).Value = "Healthcare Practitioners and Technical Occupations"
$ws.Cells.Item(13, 2).Value = "Healthcare Support Occupations"
$ws.Cells.Item(14, 2).Value = "Legal Occupations"
$ws.Cells.Item(15, 2).Value = "Community and Social Service Occupations"
$ws.Cells.Item(16, 2).Value = "Life, Physical, and Social Science Occupations"
$ws.Cells.Item(17, 2).Value = "Life, Physical, and Social Science Occupations"
$ws.Cells.Item(18, 2).Value = "Life, Physical, and Social Science Occupations"
